# Natmi following Dr Hou advice
# Rebuilds the LR-pair matrix for Tnc -> Itga7 so the sending/target
# clusters cover ECs, FAPs and sCs (previously only FAPs and sCs were
# sending clusters). Existing rows 2-7 are updated in place and three
# new rows (8-10) are appended for the ECs sending-cluster combinations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Tnc"
$ws.Cells.Item(2, 3).Value = "Itga7"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 1.442371333333333
$ws.Cells.Item(2, 8).Value = 4.327114
$ws.Cells.Item(2, 9).Value = 0.03522044016446201
$ws.Cells.Item(2, 10).Value = 0.03522044016446201
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 3.349417
$ws.Cells.Item(2, 14).Value = 10.048251
$ws.Cells.Item(2, 15).Value = 0.0602955800561437
$ws.Cells.Item(2, 16).Value = 0.0602955800561437
$ws.Cells.Item(2, 17).Value = 4.831103064179334
$ws.Cells.Item(2, 18).Value = 43.479927577614
$ws.Cells.Item(2, 19).Value = 0.002123636869548938
$ws.Cells.Item(2, 20).Value = 0.002123636869548938

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Tnc"
$ws.Cells.Item(3, 3).Value = "Itga7"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 1.442371333333333
$ws.Cells.Item(3, 8).Value = 4.327114
$ws.Cells.Item(3, 9).Value = 0.03522044016446201
$ws.Cells.Item(3, 10).Value = 0.03522044016446201
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 1.282876
$ws.Cells.Item(3, 14).Value = 3.848628
$ws.Cells.Item(3, 15).Value = 0.023094094452887
$ws.Cells.Item(3, 16).Value = 0.02309409445288699
$ws.Cells.Item(3, 17).Value = 1.850383566621333
$ws.Cells.Item(3, 18).Value = 16.653452099592
$ws.Cells.Item(3, 19).Value = 0.0008133841718303405
$ws.Cells.Item(3, 20).Value = 0.0008133841718303404

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Tnc"
$ws.Cells.Item(4, 3).Value = "Itga7"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 1.442371333333333
$ws.Cells.Item(4, 8).Value = 4.327114
$ws.Cells.Item(4, 9).Value = 0.03522044016446201
$ws.Cells.Item(4, 10).Value = 0.03522044016446201
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 50.917666
$ws.Cells.Item(4, 14).Value = 152.752998
$ws.Cells.Item(4, 15).Value = 0.9166103254909692
$ws.Cells.Item(4, 16).Value = 0.9166103254909692
$ws.Cells.Item(4, 17).Value = 73.44218179864133
$ws.Cells.Item(4, 18).Value = 660.9796361877719
$ws.Cells.Item(4, 19).Value = 0.03228341912308273
$ws.Cells.Item(4, 20).Value = 0.03228341912308273

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Tnc"
$ws.Cells.Item(5, 3).Value = "Itga7"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 11.331397
$ws.Cells.Item(5, 8).Value = 33.994191
$ws.Cells.Item(5, 9).Value = 0.2766948987373093
$ws.Cells.Item(5, 10).Value = 0.2766948987373092
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 3.349417
$ws.Cells.Item(5, 14).Value = 10.048251
$ws.Cells.Item(5, 15).Value = 0.0602955800561437
$ws.Cells.Item(5, 16).Value = 0.0602955800561437
$ws.Cells.Item(5, 17).Value = 37.95357374554901
$ws.Cells.Item(5, 18).Value = 341.582163709941
$ws.Cells.Item(5, 19).Value = 0.01668347941794201
$ws.Cells.Item(5, 20).Value = 0.016683479417942

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Tnc"
$ws.Cells.Item(6, 3).Value = "Itga7"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 11.331397
$ws.Cells.Item(6, 8).Value = 33.994191
$ws.Cells.Item(6, 9).Value = 0.2766948987373093
$ws.Cells.Item(6, 10).Value = 0.2766948987373092
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 1.282876
$ws.Cells.Item(6, 14).Value = 3.848628
$ws.Cells.Item(6, 15).Value = 0.023094094452887
$ws.Cells.Item(6, 16).Value = 0.02309409445288699
$ws.Cells.Item(6, 17).Value = 14.536777257772
$ws.Cells.Item(6, 18).Value = 130.830995319948
$ws.Cells.Item(6, 19).Value = 0.006390018126071424
$ws.Cells.Item(6, 20).Value = 0.006390018126071421

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Tnc"
$ws.Cells.Item(7, 3).Value = "Itga7"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 11.331397
$ws.Cells.Item(7, 8).Value = 33.994191
$ws.Cells.Item(7, 9).Value = 0.2766948987373093
$ws.Cells.Item(7, 10).Value = 0.2766948987373092
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 50.917666
$ws.Cells.Item(7, 14).Value = 152.752998
$ws.Cells.Item(7, 15).Value = 0.9166103254909692
$ws.Cells.Item(7, 16).Value = 0.9166103254909692
$ws.Cells.Item(7, 17).Value = 576.968287759402
$ws.Cells.Item(7, 18).Value = 5192.714589834618
$ws.Cells.Item(7, 19).Value = 0.2536214011932958
$ws.Cells.Item(7, 20).Value = 0.2536214011932958

# Row 8
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Tnc"
$ws.Cells.Item(8, 3).Value = "Itga7"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 28.17890933333333
$ws.Cells.Item(8, 8).Value = 84.536728
$ws.Cells.Item(8, 9).Value = 0.6880846610982287
$ws.Cells.Item(8, 10).Value = 0.6880846610982286
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 3.349417
$ws.Cells.Item(8, 14).Value = 10.048251
$ws.Cells.Item(8, 15).Value = 0.0602955800561437
$ws.Cells.Item(8, 16).Value = 0.0602955800561437
$ws.Cells.Item(8, 17).Value = 94.38291796252534
$ws.Cells.Item(8, 18).Value = 849.446261662728
$ws.Cells.Item(8, 19).Value = 0.04148846376865276
$ws.Cells.Item(8, 20).Value = 0.04148846376865275

# Row 9
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Tnc"
$ws.Cells.Item(9, 3).Value = "Itga7"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 28.17890933333333
$ws.Cells.Item(9, 8).Value = 84.536728
$ws.Cells.Item(9, 9).Value = 0.6880846610982287
$ws.Cells.Item(9, 10).Value = 0.6880846610982286
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 1.282876
$ws.Cells.Item(9, 14).Value = 3.848628
$ws.Cells.Item(9, 15).Value = 0.023094094452887
$ws.Cells.Item(9, 16).Value = 0.02309409445288699
$ws.Cells.Item(9, 17).Value = 36.15004648990934
$ws.Cells.Item(9, 18).Value = 325.350418409184
$ws.Cells.Item(9, 19).Value = 0.01589069215498523
$ws.Cells.Item(9, 20).Value = 0.01589069215498523

# Row 10
$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Tnc"
$ws.Cells.Item(10, 3).Value = "Itga7"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 28.17890933333333
$ws.Cells.Item(10, 8).Value = 84.536728
$ws.Cells.Item(10, 9).Value = 0.6880846610982287
$ws.Cells.Item(10, 10).Value = 0.6880846610982286
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 50.917666
$ws.Cells.Item(10, 14).Value = 152.752998
$ws.Cells.Item(10, 15).Value = 0.9166103254909692
$ws.Cells.Item(10, 16).Value = 0.9166103254909692
$ws.Cells.Item(10, 17).Value = 1434.804293678949
$ws.Cells.Item(10, 18).Value = 12913.23864311054
$ws.Cells.Item(10, 19).Value = 0.6307055051745907
$ws.Cells.Item(10, 20).Value = 0.6307055051745906
